# Daily attendance processing - 2025-10-21 11:17:38
#
# The "Recorded By" column (G) lists the people/processes that recorded
# each attendance session as a comma-separated string (e.g.
# "System, dnasr281@gmail.com"). This pass normalizes those entries by
# reversing the order of the comma-separated names in every populated
# cell of column G (header row excluded), so the most recently-added
# recorder name appears first in the underlying data and "System" is
# moved to the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Length -gt 1) {
            $reversed = @()
            for ($i = $parts.Length - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }
            $cell.Value2 = [string]::Join(", ", $reversed)
        }
    }
}
